# Auto-generated edit script applying diff changes to 广州-漫展信息.xlsx
$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 914
$ws.Range("F4").Value = 81
$ws.Range("F5").Value = 0
$ws.Range("F7").Value = 77
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = "已售罄"
$ws.Range("F10").Value = 0
$ws.Range("F12").Value = 500
$ws.Range("F13").Value = 660
$ws.Range("F14").Value = 525
$ws.Range("F17").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("F20").Value = 76
$ws.Range("F21").Value = 522
$ws.Range("F23").Value = 0
$ws.Range("F24").Value = 77
$ws.Range("F26").Value = 366
$ws.Range("F27").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("F31").Value = 345
$ws.Range("F32").Value = 0
$ws.Range("F33").Value = 797
$ws.Range("F34").Value = 347
$ws.Range("F35").Value = 166
$ws.Range("F36").Value = 200
$ws.Range("F40").Value = 948
$ws.Range("F41").Value = 0
$ws.Range("F42").Value = 68
$ws.Range("F43").Value = 0

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("F13").Value = 43
$ws.Range("F14").Value = 66
$ws.Range("F17").Value = 4351

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 152

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1682
$ws.Range("F3").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("F6").Value = 914
$ws.Range("F9").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("G11").Value = "已售罄"
$ws.Range("F17").Value = 5
$ws.Range("F19").Value = 147
$ws.Range("F20").Value = 500
$ws.Range("F22").Value = 0
$ws.Range("F24").Value = 0
$ws.Range("F26").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("F31").Value = 189
$ws.Range("F32").Value = 0
$ws.Range("F33").Value = 0
$ws.Range("F34").Value = 509
$ws.Range("F37").Value = 539
$ws.Range("F38").Value = 0
$ws.Range("F39").Value = 0
$ws.Range("F41").Value = 797
$ws.Range("F42").Value = 0
$ws.Range("F45").Value = 0
$ws.Range("F47").Value = 0
$ws.Range("F48").Value = 0
